$wb = $excel.ActiveWorkbook
$sys = $wb.Worksheets.Item("#system")

$IValues = @("assertPassword(file)", "clear(file,worksheet,range)", "clearPassword(file,password)", "columnarCsv(file,worksheet,ranges,output)", "csv(file,worksheet,range,output)", "json(file,worksheet,range,header,output)", "saveData(var,file,worksheet,range)", "saveRange(var,file,worksheet,range)", "saveTotalDataCount(file,worksheet,saveVar)", "setPassword(file,password)", "write(file,worksheet,startCell,data)", "writeAcross(file,worksheet,startCell,array)", "writeDown(file,worksheet,startCell,array)", "writeVar(var,file,worksheet,startCell)")
$SValues = @("assertContentEqual(actualPdf,expectedPdf)", "assertFormElementPresent(var,name)", "assertFormValue(var,name,expected)", "assertFormValues(var,name,expectedValues,exactOrder)", "assertPatternNotPresent(pdf,regex)", "assertPatternPresent(pdf,regex)", "assertTextArray(pdf,textArray,ordered)", "assertTextNotPresent(pdf,text)", "assertTextPresent(pdf,text)", "count(pdf,text,var)", "saveAsPages(pdf,destination)", "saveAsPdf(profile,content,file)", "saveAsText(pdf,destination)", "saveFormValues(pdf,var,pageAndLineStartEnd,strategy)", "saveMetadata(pdf,var)", "saveToVar(pdf,var)")
$ZValues = @("assertAndClick(locator,label)", "assertAttribute(locator,attrName,value)", "assertAttributeContain(locator,attrName,contains)", "assertAttributeNotContain(locator,attrName,contains)", "assertAttributeNotPresent(locator,attrName)", "assertAttributePresent(locator,attrName)", "assertChecked(locator)", "assertContainCount(locator,text,count)", "assertCssNotPresent(locator,property)", "assertCssPresent(locator,property,value)", "assertElementByAttributes(nameValues)", "assertElementByText(locator,text)", "assertElementCount(locator,count)", "assertElementEnabled(locator)", "assertElementNotPresent(locator)", "assertElementPresent(locator)", "assertElementsPresent(prefix)", "assertFocus(locator)", "assertFrameCount(count)", "assertFramePresent(frameName)", "assertIECompatMode()", "assertIENativeMode()", "assertLinkByLabel(label)", "assertMultiSelect(locator)", "assertNotChecked(locator)", "assertNotFocus(locator)", "assertNotText(locator,text)", "assertNotVisible(locator)", "assertOneMatch(locator)", "assertScrollbarHNotPresent(locator)", "assertScrollbarHPresent(locator)", "assertScrollbarVNotPresent(locator)", "assertScrollbarVPresent(locator)", "assertSingleSelect(locator)", "assertTable(locator,row,column,text)", "assertText(locator,text)", "assertTextContains(locator,text)", "assertTextCount(locator,text,count)", "assertTextList(locator,list,ignoreOrder)", "assertTextMatches(text,minMatch,scrollTo)", "assertTextNotContain(locator,text)", "assertTextNotPresent(text)", "assertTextOrder(locator,descending)", "assertTextPresent(text)", "assertTitle(text)", "assertValue(locator,value)", "assertValueOrder(locator,descending)", "assertVisible(locator)", "checkAll(locator)", "clearLocalStorage()", "click(locator)", "clickAll(locator)", "clickAndWait(locator,waitMs)", "clickByLabel(label)", "clickByLabelAndWait(label,waitMs)", "clickOffset(locator,x,y)", "clickWithKeys(locator,keys)", "close()", "closeAll()", "deselect(locator,text)", "deselectMulti(locator,array)", "dismissInvalidCert()", "dismissInvalidCertPopup()", "doubleClick(locator)", "doubleClickAndWait(locator,waitMs)", "doubleClickByLabel(label)", "doubleClickByLabelAndWait(label,waitMs)", "dragAndDrop(fromLocator,toLocator)", "dragTo(fromLocator,xOffset,yOffset)", "editLocalStorage(key,value)", "executeScript(var,script)", "focus(locator)", "goBack()", "goBackAndWait()", "maximizeWindow()", "mouseOver(locator)", "open(url)", "openAndWait(url,waitMs)", "openHttpBasic(url,username,password)", "openIgnoreTimeout(url)", "refresh()", "refreshAndWait()", "resizeWindow(width,height)", "rightClick(locator)", "saveAllWindowIds(var)", "saveAllWindowNames(var)", "saveAttribute(var,locator,attrName)", "saveAttributeList(var,locator,attrName)", "saveBrowserVersion(var)", "saveCount(var,locator)", "saveDivsAsCsv(headers,rows,cells,nextPage,file)", "saveElement(var,locator)", "saveElements(var,locator)", "saveInfiniteDivsAsCsv(config,file)", "saveInfiniteTableAsCsv(config,file)", "saveLocalStorage(var,key)", "saveLocation(var)", "savePageAs(var,sessionIdName,url)", "savePageAsFile(sessionIdName,url,file)", "saveSelectedText(var,locator)", "saveSelectedValue(var,locator)", "saveTableAsCsv(locator,nextPageLocator,file)", "saveText(var,locator)", "saveTextArray(var,locator)", "saveTextSubstringAfter(var,locator,delim)", "saveTextSubstringBefore(var,locator,delim)", "saveTextSubstringBetween(var,locator,start,end)", "saveTitle(var)", "saveValue(var,locator)", "saveValues(var,locator)", "screenshot(file,locator)", "scrollElement(locator,xOffset,yOffset)", "scrollLeft(locator,pixel)", "scrollPage(xOffset,yOffset)", "scrollRight(locator,pixel)", "scrollTo(locator)", "select(locator,text)", "selectAllOptions(locator)", "selectFrame(locator)", "selectMulti(locator,array)", "selectMultiByValue(locator,array)", "selectMultiOptions(locator)", "selectText(locator)", "selectWindow(winId)", "selectWindowAndWait(winId,waitMs)", "selectWindowByIndex(index)", "selectWindowByIndexAndWait(index,waitMs)", "switchBrowser(profile,config)", "toggleSelections(locator)", "type(locator,value)", "typeKeys(locator,value)", "uncheckAll(locator)", "unselectAllText()", "updateAttribute(locator,attrName,value)", "upload(fieldLocator,file)", "verifyContainText(locator,text)", "verifyText(locator,text)", "wait(waitMs)", "waitForElementPresent(locator)", "waitForElementsPresent(locators)", "waitForPopUp(winId,waitMs)", "waitForTextPresent(text)", "waitForTitle(text)")
$ABValues = @("assertNotPresent(name)", "assertPresent(name)", "assertValue(name,value)", "clearCookieFields(var,remove)", "delete(name)", "deleteAll()", "save(var,name)", "saveAll(var)", "saveAllAsText(var,exclude)")

$colMap = @{ "I" = 9; "S" = 19; "Z" = 26; "AB" = 28 }

function Set-ColumnValues($colLetter, $values) {
    $colIndex = $colMap[$colLetter]
    for ($i = 0; $i -lt $values.Count; $i++) {
        $sys.Cells.Item($i + 2, $colIndex).Value2 = $values[$i]
    }
}

Set-ColumnValues "I" $IValues
Set-ColumnValues "S" $SValues
Set-ColumnValues "Z" $ZValues
Set-ColumnValues "AB" $ABValues

# Update named ranges to reflect the new extents
$wb.Names.Item("excel").RefersTo = "='#system'!`$I`$2:`$I`$15"
$wb.Names.Item("pdf").RefersTo = "='#system'!`$S`$2:`$S`$17"
$wb.Names.Item("web").RefersTo = "='#system'!`$Z`$2:`$Z`$144"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AB`$2:`$AB`$10"
